$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above row 37 for the latest weekly data point; ---
# --- this pushes the existing rows 37-50 down to become rows 38-51.  ---
$ws.Rows.Item(37).Insert()

# --- Populate the newly inserted row 37 ---
$ws.Range("A37").Value = 6
$ws.Range("B37").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C37").Value = "Metropolitana"
$ws.Range("D37").Value = 44943
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100102
$ws.Range("H37").Value = "Cítricos"
$ws.Range("I37").Value = 100102006
$ws.Range("J37").Value = "Pomelo"
$ws.Range("K37").Value = "Start Ruby"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 24
$ws.Range("N37").Value = 180000
$ws.Range("O37").Value = 200000
$ws.Range("P37").Value = 190000
$ws.Range("Q37").Value = "$/bins (350 kilos)"
$ws.Range("R37").Value = "Provincia de Limarí"
$ws.Range("S37").Value = 543
$ws.Range("T37").Value = 350
